# Update NATMI LR-pair output (Fasl-Fas) with newly recomputed TPM-based values.
# Sending cluster changes from "ECs" to "MuSCs" for every row, while the
# Target cluster values are re-assigned (ECs, FAPs, MuSCs) and all of the
# derived expression / specificity metrics are refreshed with new numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 : MuSCs -> Fasl -> Fas -> ECs -------------------------------
$ws.Range("A2").Value = "MuSCs"
$ws.Range("B2").Value = "Fasl"
$ws.Range("C2").Value = "Fas"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.01577866666666667
$ws.Range("H2").Value = 0.047336
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 19.58374133333333
$ws.Range("N2").Value = 58.751224
$ws.Range("O2").Value = 0.6578841466750758
$ws.Range("P2").Value = 0.6578841466750758
$ws.Range("Q2").Value = 0.3090053265848889
$ws.Range("R2").Value = 2.781047939264
$ws.Range("S2").Value = 0.6578841466750758
$ws.Range("T2").Value = 0.6578841466750758

# --- Row 3 : MuSCs -> Fasl -> Fas -> FAPs -------------------------------
$ws.Range("A3").Value = "MuSCs"
$ws.Range("B3").Value = "Fasl"
$ws.Range("C3").Value = "Fas"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.01577866666666667
$ws.Range("H3").Value = 0.047336
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 6.657374333333333
$ws.Range("N3").Value = 19.972123
$ws.Range("O3").Value = 0.2236437337398222
$ws.Range("P3").Value = 0.2236437337398222
$ws.Range("Q3").Value = 0.1050444904808889
$ws.Range("R3").Value = 0.945400414328
$ws.Range("S3").Value = 0.2236437337398222
$ws.Range("T3").Value = 0.2236437337398222

# --- Row 4 : MuSCs -> Fasl -> Fas -> MuSCs ------------------------------
$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Fasl"
$ws.Range("C4").Value = "Fas"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.01577866666666667
$ws.Range("H4").Value = 0.047336
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.526650333333333
$ws.Range("N4").Value = 10.579951
$ws.Range("O4").Value = 0.118472119585102
$ws.Range("P4").Value = 0.118472119585102
$ws.Range("Q4").Value = 0.05564584005955555
$ws.Range("R4").Value = 0.500812560536
$ws.Range("S4").Value = 0.118472119585102
$ws.Range("T4").Value = 0.118472119585102
